$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: Katie Kupferschmid - move from 3rd Interview to 4th Interview, update date
$ws.Range("E5").Value = "4th Interview"
$ws.Range("F5").Value = 45992

# Row 8 previously held Mads Berli (CV Sent); he is removed and the prior
# row 9 (Zemir  Sadikovic) moves up into row 8, with status advancing to
# 1st Interview and the date updated.
$ws.Range("D8").Value = "Zemir  Sadikovic"
$ws.Range("E8").Value = "1st Interview"
$ws.Range("F8").Value = 45992

# Delete the now-duplicate former row 9, shifting remaining rows up
$ws.Rows(9).Delete()
